# "Add remarks to grade form"
# Fill in / update the student-remarks column (H) on the Evaluation form sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluation form")

# Novice skills section (rows 22-29) - column H holds the teacher's remarks
# next to each indicator.
$ws.Range("H22").Value = "Tree(spheres, cylinders), car, house, playground,plane and cone bush"
$ws.Range("H23").Value = "Spheres, cylinders, plane, cone"
$ws.Range("H26").Value = "plane, tree leaves, tree log, bush leaves"
$ws.Range("H29").Value = "Orbit, flying and first person"

# Intermediate / expert skills section (rows 34-35)
$ws.Range("H34").Value = "comments, readable, multiple files, modular, multiple functions and a clear structure"
$ws.Range("H35").Value = "Quaint little cabin in the woods with a moving car and a playground"

# Move the view back to the top of the sheet and reset the active selection
# (matches the author re-opening the form at the total-score cell).
$ws.Range("H19").Select() | Out-Null
